# "add thumb field in bulk upload"
#
# Adds a new "thumbnail_id" column (P) to the product bulk-upload demo
# sheet: a styled header in P1 and a sample numeric value in P2, matching
# the look of the other "special" headers already on the sheet (D1, E1,
# M1, O1 - bold Arial, dark color), plus left alignment for the new
# header specifically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell: "thumbnail_id"
$ws.Range("P1").Value = "thumbnail_id"

# Match the header styling used by the other special header cells
# (D1/E1/M1/O1: bold Arial, dark text) by copying D1's formatting onto
# the new header cell...
$ws.Range("D1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats

# ...then left-align this particular header.
$ws.Range("P1").HorizontalAlignment = -4131   # xlLeft

# Sample data value for the new column, on the existing demo data row.
$ws.Range("P2").Value = 142

# Keep the active selection where the saved workbook leaves it.
$ws.Range("D3").Select()
